# Auto-generated edit script: apply cell value changes per diff
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1174.5714
$ws.Range("I20").Value = 537
$ws.Range("K20").Value = 537
$ws.Range("M20").Value = -307
$ws.Range("H35").Value = 1174.5714
$ws.Range("I35").Value = 537
$ws.Range("K35").Value = 537
$ws.Range("M35").Value = -158
$ws.Range("H51").Value = 14999.5
$ws.Range("I51").Value = 14999.5
$ws.Range("K51").Value = 14999.5
$ws.Range("M51").Value = -14515.5
$ws.Range("H62").Value = 14222
$ws.Range("I62").Value = 14000
$ws.Range("K62").Value = 14000
$ws.Range("M62").Value = -13376
$ws.Range("H65").Value = 14222
$ws.Range("I65").Value = 14000
$ws.Range("K65").Value = 70000
$ws.Range("M65").Value = -66880
$ws.Range("H68").Value = 65000
$ws.Range("J68").Value = 65000
$ws.Range("L68").Value = 65000
$ws.Range("N68").Value = -66498
$ws.Range("H70").Value = 1238.625
$ws.Range("I70").Value = 1300
$ws.Range("J70").Value = 1218.1666
$ws.Range("K70").Value = 3900
$ws.Range("L70").Value = 3654.4998
$ws.Range("M70").Value = -3630
$ws.Range("N70").Value = -4194.4998
$ws.Range("H71").Value = 65000
$ws.Range("J71").Value = 65000
$ws.Range("L71").Value = 195000
$ws.Range("N71").Value = -202488
$ws.Range("H73").Value = 1238.625
$ws.Range("I73").Value = 1300
$ws.Range("J73").Value = 1218.1666
$ws.Range("K73").Value = 3900
$ws.Range("L73").Value = 3654.4998
$ws.Range("M73").Value = -2964
$ws.Range("N73").Value = -5526.4998
$ws.Range("H86").Value = 10334.333
$ws.Range("I86").Value = 10501.5
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 10501.5
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -9378.5
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 10334.333
$ws.Range("I89").Value = 10501.5
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 52507.5
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -46891.5
$ws.Range("N89").Value = -61232
$ws.Range("H116").Value = 6669.1924
$ws.Range("I116").Value = 6445.364
$ws.Range("K116").Value = 6445.364
$ws.Range("M116").Value = -3003.364
$ws.Range("H125").Value = 83339680
$ws.Range("I125").Value = 250000000
$ws.Range("K125").Value = 2250000000
$ws.Range("M125").Value = -2249997540
$ws.Range("H135").Value = 1142.4286
$ws.Range("I135").Value = 1149.5
$ws.Range("K135").Value = 10345.5
$ws.Range("M135").Value = -7810.5
$ws.Range("H137").Value = 3785.5715
$ws.Range("J137").Value = 3785.5715
$ws.Range("L137").Value = 11356.7145
$ws.Range("N137").Value = -16456.7145
$ws.Range("H138").Value = 4999.225
$ws.Range("J138").Value = 5679.1
$ws.Range("L138").Value = 17037.3
$ws.Range("N138").Value = -27317.3

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 462.75
$ws.Range("I6").Value = 462.75
$ws.Range("K6").Value = 462.75
$ws.Range("M6").Value = -289.75
$ws.Range("H19").Value = 1800
$ws.Range("I19").Value = 1500
$ws.Range("J19").Value = 1950
$ws.Range("K19").Value = 1500
$ws.Range("L19").Value = 1950
$ws.Range("M19").Value = -1271
$ws.Range("N19").Value = -2408
$ws.Range("H26").Value = 321.4
$ws.Range("I26").Value = 321.4
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 321.4
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 8.600000000000023
$ws.Range("N26").ClearContents()
$ws.Range("H37").Value = 20508.5
$ws.Range("I37").Value = 20508.5
$ws.Range("K37").Value = 20508.5
$ws.Range("M37").Value = -20235.5
$ws.Range("H46").Value = 23500
$ws.Range("I46").Value = 25000
$ws.Range("J46").Value = 22000
$ws.Range("K46").Value = 25000
$ws.Range("L46").Value = 22000
$ws.Range("M46").Value = -24681
$ws.Range("N46").Value = -22638
$ws.Range("H61").Value = 2349.5
$ws.Range("I61").Value = 2349.5
$ws.Range("K61").Value = 2349.5
$ws.Range("M61").Value = -2137.5
$ws.Range("H64").Value = 54999
$ws.Range("J64").Value = 54999
$ws.Range("L64").Value = 54999
$ws.Range("N64").Value = -55495
$ws.Range("H67").Value = 54999
$ws.Range("J67").Value = 54999
$ws.Range("L67").Value = 54999
$ws.Range("N67").Value = -56715
$ws.Range("H74").Value = 1162.8572
$ws.Range("I74").Value = 1106.6666
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1106.6666
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -232.6666
$ws.Range("N74").Value = -3248
$ws.Range("H77").Value = 1162.8572
$ws.Range("I77").Value = 1106.6666
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 5533.333000000001
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -1165.333000000001
$ws.Range("N77").Value = -16236
$ws.Range("H132").Value = 3201.4546
$ws.Range("I132").Value = 2062.75
$ws.Range("K132").Value = 6188.25
$ws.Range("M132").Value = -3658.25
$ws.Range("H136").Value = 2349.5
$ws.Range("I136").Value = 2349.5
$ws.Range("K136").Value = 7048.5
$ws.Range("M136").Value = -4498.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("K6").Value = 500
$ws.Range("M6").Value = -387
$ws.Range("H16").Value = 1023.125
$ws.Range("I16").Value = 875.75
$ws.Range("K16").Value = 875.75
$ws.Range("M16").Value = -588.75
$ws.Range("H18").Value = 79965.7
$ws.Range("J18").Value = 79965.7
$ws.Range("L18").Value = 79965.7
$ws.Range("N18").Value = -80425.7
$ws.Range("H56").Value = 49998.5
$ws.Range("J56").Value = 49998.5
$ws.Range("L56").Value = 49998.5
$ws.Range("N56").Value = -51688.5
$ws.Range("H58").Value = 15007
$ws.Range("J58").Value = 15007
$ws.Range("L58").Value = 15007
$ws.Range("N58").Value = -15413
$ws.Range("H113").Value = 1023.125
$ws.Range("I113").Value = 875.75
$ws.Range("K113").Value = 875.75
$ws.Range("M113").Value = 1294.25
$ws.Range("H136").Value = 15007
$ws.Range("J136").Value = 15007
$ws.Range("L136").Value = 45021
$ws.Range("N136").Value = -50121

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 3833.3333
$ws.Range("I35").Value = 2000
$ws.Range("J35").Value = 4357.143
$ws.Range("K35").Value = 6000
$ws.Range("L35").Value = 13071.429
$ws.Range("M35").Value = -5712
$ws.Range("N35").Value = -13647.429
$ws.Range("H136").Value = 7325.3335
$ws.Range("I136").Value = 7325.3335
$ws.Range("K136").Value = 21976.0005
$ws.Range("M136").Value = -16876.0005

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 663.4
$ws.Range("I102").Value = 579.25
$ws.Range("K102").Value = 579.25
$ws.Range("M102").Value = 1042.75
$ws.Range("H122").Value = 201091.4
$ws.Range("I122").Value = 251064.25
$ws.Range("K122").Value = 753192.75
$ws.Range("M122").Value = -750742.75

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 453.85715
$ws.Range("I30").Value = 453.85715
$ws.Range("K30").Value = 453.85715
$ws.Range("M30").Value = -345.85715
$ws.Range("H40").Value = 9948.5
$ws.Range("I40").Value = 5899.5
$ws.Range("K40").Value = 5899.5
$ws.Range("M40").Value = -5763.5
$ws.Range("H46").Value = 2271.6216
$ws.Range("I46").Value = 2257.1428
$ws.Range("J46").Value = 2275
$ws.Range("K46").Value = 2257.1428
$ws.Range("L46").Value = 2275
$ws.Range("M46").Value = -2069.1428
$ws.Range("N46").Value = -2651
$ws.Range("H55").Value = 963.5
$ws.Range("I55").Value = 216.8
$ws.Range("J55").Value = 1710.2
$ws.Range("K55").Value = 216.8
$ws.Range("L55").Value = 1710.2
$ws.Range("M55").Value = -43.80000000000001
$ws.Range("N55").Value = -2056.2
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H104").Value = 60000
$ws.Range("J104").Value = 60000
$ws.Range("L104").Value = 60000
$ws.Range("N104").Value = -66988
$ws.Range("H135").Value = 60999
$ws.Range("J135").Value = 60999
$ws.Range("L135").Value = 60999
$ws.Range("N135").Value = -71139

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 19999
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H62").Value = 4867.6665
$ws.Range("J62").Value = 6151.5
$ws.Range("L62").Value = 6151.5
$ws.Range("N62").Value = -7399.5
$ws.Range("H65").Value = 4867.6665
$ws.Range("J65").Value = 6151.5
$ws.Range("L65").Value = 30757.5
$ws.Range("N65").Value = -36997.5
$ws.Range("H81").Value = 1098
$ws.Range("I81").Value = 1098
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2196
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1135
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1098
$ws.Range("I84").Value = 1098
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 1098
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -5676
$ws.Range("N84").ClearContents()
